# Fix the double minus on the third derivate relation
# Negate every value in column F for rows 2 through 479 (flip the sign).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 479; $r++) {
    $cell = $ws.Cells.Item($r, 6)   # column F is the 6th column
    $val = $cell.Value2
    if ($null -ne $val) {
        $cell.Value2 = -1 * $val
    }
}
